$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1536
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H19").Value = 438.4
$ws.Range("I19").Value = 393.1
$ws.Range("J19").Value = 483.7
$ws.Range("K19").Value = 393.1
$ws.Range("L19").Value = 483.7
$ws.Range("M19").Value = -218.1
$ws.Range("N19").Value = -833.7

$ws.Range("H40").Value = 4169.8335
$ws.Range("J40").Value = 4369
$ws.Range("L40").Value = 4369
$ws.Range("N40").Value = -4719

$ws.Range("H86").Value = 3456.3333
$ws.Range("I86").Value = 2333.3333
$ws.Range("K86").Value = 2333.3333
$ws.Range("M86").Value = -1210.3333

$ws.Range("H89").Value = 3456.3333
$ws.Range("I89").Value = 2333.3333
$ws.Range("K89").Value = 11666.6665
$ws.Range("M89").Value = -6050.666499999999

$ws.Range("H96").Value = 656.2941
$ws.Range("J96").Value = 1028.2858
$ws.Range("L96").Value = 3084.8574
$ws.Range("N96").Value = -5830.857400000001

$ws.Range("H112").Value = 1811.9333
$ws.Range("I112").Value = 1575
$ws.Range("J112").Value = 1871.1666
$ws.Range("K112").Value = 4725
$ws.Range("L112").Value = 5613.4998
$ws.Range("M112").Value = -3617
$ws.Range("N112").Value = -7829.4998

$ws.Range("H116").Value = 19071.715
$ws.Range("I116").Value = 17749.75
$ws.Range("J116").Value = 19600.5
$ws.Range("K116").Value = 17749.75
$ws.Range("L116").Value = 19600.5
$ws.Range("M116").Value = -14307.75
$ws.Range("N116").Value = -26484.5

$ws.Range("H136").Value = 349969.12
$ws.Range("J136").Value = 349969.12
$ws.Range("L136").Value = 349969.12
$ws.Range("N136").Value = -360169.12

$ws.Range("H138").Value = 2186209.2
$ws.Range("I138").Value = 15628.444
$ws.Range("J138").Value = 2651333.5
$ws.Range("K138").Value = 46885.33199999999
$ws.Range("L138").Value = 7954000.5
$ws.Range("M138").Value = -41745.33199999999
$ws.Range("N138").Value = -7964280.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 172079.39
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

$ws.Range("H74").Value = 6120.1665
$ws.Range("I74").Value = 5744.2
$ws.Range("K74").Value = 5744.2
$ws.Range("M74").Value = -4870.2

$ws.Range("H77").Value = 6120.1665
$ws.Range("I77").Value = 5744.2
$ws.Range("K77").Value = 28721
$ws.Range("M77").Value = -24353

$ws.Range("H97").Value = 495.97562
$ws.Range("I97").Value = 399.0606
$ws.Range("K97").Value = 399.0606
$ws.Range("M97").Value = 96.93939999999998

$ws.Range("H102").Value = 6623.654
$ws.Range("I102").Value = 5189.7646
$ws.Range("J102").Value = 9332.111000000001
$ws.Range("K102").Value = 5189.7646
$ws.Range("L102").Value = 9332.111000000001
$ws.Range("M102").Value = -3567.7646
$ws.Range("N102").Value = -12576.111

$ws.Range("H122").Value = 9261767
$ws.Range("I122").Value = 13891201
$ws.Range("K122").Value = 41673603
$ws.Range("M122").Value = -41671153

$ws.Range("H139").Value = 246037.25
$ws.Range("J139").Value = 246037.25
$ws.Range("L139").Value = 246037.25
$ws.Range("N139").Value = -256317.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 30667.666
$ws.Range("J96").Value = 79995
$ws.Range("L96").Value = 79995
$ws.Range("N96").Value = -85487

$ws.Range("H134").Value = 4163.6665
$ws.Range("I134").Value = 3788.182
$ws.Range("J134").Value = 4753.7144
$ws.Range("K134").Value = 11364.546
$ws.Range("L134").Value = 14261.1432
$ws.Range("M134").Value = -8829.545999999998
$ws.Range("N134").Value = -19331.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5021.1763
$ws.Range("I31").Value = 2137.5715
$ws.Range("K31").Value = 2137.5715
$ws.Range("M31").Value = -1842.5715

$ws.Range("H34").Value = 5021.1763
$ws.Range("I34").Value = 2137.5715
$ws.Range("K34").Value = 2137.5715
$ws.Range("M34").Value = -1935.5715

$ws.Range("H99").Value = 2491.5264
$ws.Range("I99").Value = 2463.2778
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 2463.2778
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -965.2777999999998
$ws.Range("N99").Value = -5996

$ws.Range("H105").Value = 4452.4443
$ws.Range("I105").Value = 5408.625
$ws.Range("K105").Value = 5408.625
$ws.Range("M105").Value = -3661.625

$ws.Range("H122").Value = 2751.889
$ws.Range("J122").Value = 3733.25
$ws.Range("L122").Value = 11199.75
$ws.Range("N122").Value = -16099.75

$ws.Range("H126").Value = 2491.5264
$ws.Range("I126").Value = 2463.2778
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 7389.8334
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -4919.8334
$ws.Range("N126").Value = -13940

$ws.Range("H141").Value = 813330
$ws.Range("J141").Value = 813330
$ws.Range("L141").Value = 813330
$ws.Range("N141").Value = -823690

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 61606816
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H14").Value = 128.26666
$ws.Range("I14").Value = 128.26666
$ws.Range("K14").Value = 384.79998
$ws.Range("M14").Value = -211.79998

$ws.Range("H86").Value = 543.5
$ws.Range("I86").Value = 175
$ws.Range("J86").Value = 666.3333
$ws.Range("K86").Value = 525
$ws.Range("L86").Value = 1998.9999
$ws.Range("M86").Value = 661
$ws.Range("N86").Value = -4370.9999

$ws.Range("H89").Value = 543.5
$ws.Range("I89").Value = 175
$ws.Range("J89").Value = 666.3333
$ws.Range("K89").Value = 1575
$ws.Range("L89").Value = 5996.9997
$ws.Range("M89").Value = 4353
$ws.Range("N89").Value = -17852.9997

$ws.Range("H107").Value = 1634.9524
$ws.Range("J107").Value = 1392.875
$ws.Range("L107").Value = 4178.625
$ws.Range("N107").Value = -8018.625

$ws.Range("H113").Value = 9654.777
$ws.Range("I113").Value = 15580
$ws.Range("J113").Value = 2248.25
$ws.Range("K113").Value = 46740
$ws.Range("L113").Value = 6744.75
$ws.Range("M113").Value = -44570
$ws.Range("N113").Value = -11084.75

$ws.Range("H117").Value = 2451.4443
$ws.Range("I117").Value = 2812.6
$ws.Range("K117").Value = 8437.799999999999
$ws.Range("M117").Value = -4995.799999999999

$ws.Range("H121").Value = 4798.5713
$ws.Range("I121").Value = 572.5
$ws.Range("J121").Value = 10433.333
$ws.Range("K121").Value = 1717.5
$ws.Range("L121").Value = 31299.999
$ws.Range("M121").Value = -407.5
$ws.Range("N121").Value = -33919.999

$ws.Range("H127").Value = 2795.1667
$ws.Range("J127").Value = 2795.1667
$ws.Range("L127").Value = 8385.500100000001
$ws.Range("N127").Value = -18305.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8312.546
$ws.Range("I80").Value = 4721
$ws.Range("K80").Value = 4721
$ws.Range("M80").Value = -3723

$ws.Range("H83").Value = 8312.546
$ws.Range("I83").Value = 4721
$ws.Range("K83").Value = 23605
$ws.Range("M83").Value = -18613

$ws.Range("H97").Value = 787.5263
$ws.Range("I97").Value = 762.2
$ws.Range("J97").Value = 882.5
$ws.Range("K97").Value = 762.2
$ws.Range("L97").Value = 882.5
$ws.Range("M97").Value = -266.2
$ws.Range("N97").Value = -1874.5

$ws.Range("H132").Value = 2984.1428
$ws.Range("I132").Value = 2737.9333
$ws.Range("K132").Value = 8213.7999
$ws.Range("M132").Value = -5683.7999

$ws.Range("H140").Value = 348311.8
$ws.Range("J140").Value = 348311.8
$ws.Range("L140").Value = 348311.8
$ws.Range("N140").Value = -358671.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5294.1035
$ws.Range("I40").Value = 3824.5715
$ws.Range("K40").Value = 3824.5715
$ws.Range("M40").Value = -3688.5715

$ws.Range("H93").Value = 7750.4
$ws.Range("I93").Value = 3500
$ws.Range("J93").Value = 8222.666999999999
$ws.Range("K93").Value = 3500
$ws.Range("L93").Value = 8222.666999999999
$ws.Range("M93").Value = -2252
$ws.Range("N93").Value = -10718.667

$ws.Range("H122").Value = 8291.929
$ws.Range("I122").Value = 8160.5386
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 24481.6158
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -22031.6158
$ws.Range("N122").Value = -34900

$ws.Range("H132").Value = 2109014
$ws.Range("I132").Value = 2938.0833
$ws.Range("K132").Value = 8814.249899999999
$ws.Range("M132").Value = -6284.249899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3290.7856
$ws.Range("I96").Value = 1995.6666
$ws.Range("J96").Value = 3644
$ws.Range("K96").Value = 1995.6666
$ws.Range("L96").Value = 3644
$ws.Range("M96").Value = -622.6666
$ws.Range("N96").Value = -6390

$ws.Range("H132").Value = 6246.8945
$ws.Range("I132").Value = 2805.0881
$ws.Range("K132").Value = 8415.264299999999
$ws.Range("M132").Value = -5885.264299999999

$ws.Range("H140").Value = 87143
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 87143
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 87143
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -97503
